$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F values (rows 2-44): 0.5 -> -0.5
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 6).Value = -0.5
}

# Update column E values (rows 45-67): 1E-3 -> 0.05
for ($r = 45; $r -le 67; $r++) {
    $ws.Cells.Item($r, 5).Value = 0.05
}

# Update the view: scroll so topLeftCell = A34, select F46
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F46").Select()
